# Applies the 2023-05-24 GitHub Actions crypto price refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) mixes plain numbers and strings that look like
# numbers (e.g. "26.913.22"). The source file stores every one of these as
# literal text, so force the whole column range to Text before writing the
# new values - this stops Excel from "helpfully" reinterpreting them as
# numbers/dates. ClearFormats afterwards drops the Text number-format again
# so the cells keep their original (unstyled) appearance.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.913.22"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "1.833.89"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("D4").Value = "1.005"

$ws.Range("D5").Value = "310.34"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -1.43%  "

$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -1.65%  "

$ws.Range("D9").Value = "0.07172"
$ws.Range("E9").Value = "  -2.79%  "

$ws.Range("D10").Value = "0.8797"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").Value = "0.07839"
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").Value = "19.64"
$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").Value = "1.838.63"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").Value = "5.343"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "6.396"
$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").Value = "88.40"
$ws.Range("E16").Value = "  -4.62%  "

$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "0.000008749"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "26.940.38"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "5.019"
$ws.Range("E22").Value = "  -2.68%  "

$ws.Range("D23").Value = "10.44"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("D24").Value = "1.978"
$ws.Range("E24").Value = "  +5.21%  "

$ws.Range("D25").Value = "150.76"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("D26").Value = "18.18"
$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("D27").Value = "1.992"
$ws.Range("E27").Value = "  -4.54%  "

$ws.Range("D28").Value = "113.74"
$ws.Range("E28").Value = "  -2.77%  "

$ws.Range("D29").Value = "4.965"
$ws.Range("E29").Value = "  -3.89%  "

$ws.Range("D30").Value = "0.08841"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").Value = "3.125"
$ws.Range("E31").Value = "  +3.23%  "

$ws.Range("D32").Value = "0.7654"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").Value = "4.466"

$ws.Range("D34").Value = "1.133"
$ws.Range("E34").Value = "  -2.64%  "

$ws.Range("D35").Value = "2.671"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("D37").Value = "0.01938"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").Value = "2.933"
$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("D39").Value = "0.05136"
$ws.Range("E39").Value = "  -2.72%  "

$ws.Range("D40").Value = "6.944"
$ws.Range("E40").Value = "  -3.27%  "

$ws.Range("D41").Value = "0.4978"
$ws.Range("E41").Value = "  -4.49%  "

$ws.Range("D42").Value = "0.1601"
$ws.Range("E42").Value = "  -2.75%  "

$ws.Range("D43").Value = "8.313"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "0.4698"
$ws.Range("E44").Value = "  -3.59%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.004"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("D47").Value = "103.00"

$ws.Range("D48").Value = "1.615"
$ws.Range("E48").Value = "  -2.58%  "

$ws.Range("D49").Value = "0.06098"
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("D50").Value = "64.95"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("D51").Value = "36.44"
$ws.Range("E51").Value = "  -1.84%  "

$dRange.ClearFormats()
